$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (M1_PH, CM2_PH, CMN3_PH, CMN4_PH) previously blank inline-string cells,
# now populated with computed numeric values.
$ws.Range("I2").Value = -0.1510628729506125
$ws.Range("J2").Value = 0.2917876456392605
$ws.Range("K2").Value = -0.250788718282305
$ws.Range("L2").Value = 2.387835162127833

# Row 19 (M1_PH, CM2_PH, CMN3_PH, CMN4_PH) previously blank inline-string cells,
# now populated with computed numeric values.
$ws.Range("I19").Value = -0.3738853921488089
$ws.Range("J19").Value = 0.390390469269369
$ws.Range("K19").Value = 0.06626720357243029
$ws.Range("L19").Value = 1.919622844000926
